# Updates the cryptos list (prices / 1h volume %) on Sheet1, including a
# couple of rank swaps (Avalanche<->Uniswap at rows 21/22, and
# TrustWalletToken<->InternetComputer(DFINITY) at rows 37/38).
# Numeric-looking "Price" values must stay plain text, so cells whose new
# value would otherwise auto-parse as a number are briefly switched to a
# text number format, written, then restored to the default "Normal"
# style so no stray per-cell formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '20.258.70'
$ws.Cells.Item(2, 5).Value = '  +2.01%  '
$ws.Cells.Item(3, 4).Value = '1.444.91'
$ws.Cells.Item(3, 5).Value = '  +3.91%  '
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.06%  '
$ws.Cells.Item(5, 5).Value = '  -9.02%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '278.36'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +3.97%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3661'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.38%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3126'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.91%  '
$ws.Cells.Item(9, 5).Value = '  +0.34%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.024'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +7.39%  '
$ws.Cells.Item(11, 5).Value = '  +3.02%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9996'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -0.61%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.404'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +4.14%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.70'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +9.19%  '
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.076'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +1.26%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00001019'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +3.95%  '
$ws.Cells.Item(17, 4).Value = '1.444.04'
$ws.Cells.Item(17, 5).Value = '  +3.21%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9425'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -6.04%  '
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05643'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.80%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '68.51'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.77%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.415'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -0.33%  '
$ws.Cells.Item(22, 2).Value = 'Avalanche'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '14.46'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +2.81%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.88'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +4.38%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.256'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.33%  '
$ws.Cells.Item(25, 4).Value = '20.285.59'
$ws.Cells.Item(25, 5).Value = '  +2.01%  '
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.175'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.48%  '
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '137.92'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +2.11%  '
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.02'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +4.18%  '
$ws.Cells.Item(29, 4).Value = '1.597.97'
$ws.Cells.Item(29, 5).Value = '  +2.66%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '110.14'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +3.86%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.827'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.78%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8061'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +4.03%  '
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.854'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -5.57%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07702'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.42%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05954'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +7.49%  '
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.453'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +11.59%  '
$ws.Cells.Item(37, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.694'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +1.46%  '
$ws.Cells.Item(38, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.143'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +12.50%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.02001'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +1.05%  '
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.19'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +4.06%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9307'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -7.11%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1844'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.94%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.186'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -12.51%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.528'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +2.20%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5247'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +2.92%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.09'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +2.58%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '119.32'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +11.55%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5153'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +4.93%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.763'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +4.34%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06340'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +4.44%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9926'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.12%  '
